$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 13573.4
$ws.Range("I47").Value = 13573.4
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 13573.4
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -12601.4
$ws.Range("H74").Value = 4805.778
$ws.Range("I74").Value = 4208.6665
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 4208.6665
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -3272.6665
$ws.Range("N74").Value = -7872
$ws.Range("H77").Value = 4805.778
$ws.Range("I77").Value = 4208.6665
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 21043.3325
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -16363.3325
$ws.Range("N77").Value = -39360
$ws.Range("H98").Value = 3282.862
$ws.Range("I98").Value = 2850.2307
$ws.Range("K98").Value = 2850.2307
$ws.Range("M98").Value = -1352.2307
$ws.Range("H112").Value = 2367.45
$ws.Range("I112").Value = 999.8570999999999
$ws.Range("K112").Value = 2999.5713
$ws.Range("M112").Value = -1891.5713
$ws.Range("H122").Value = 3282.862
$ws.Range("I122").Value = 2850.2307
$ws.Range("K122").Value = 8550.6921
$ws.Range("M122").Value = -6100.6921
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").ClearContents()
$ws.Range("N124").Value = 0
$ws.Range("H132").Value = 4758.625
$ws.Range("I132").Value = 4839.4346
$ws.Range("K132").Value = 14518.3038
$ws.Range("M132").Value = -11988.3038
$ws.Range("H135").Value = 1030.7273
$ws.Range("I135").Value = 1069.2222
$ws.Range("K135").Value = 9622.9998
$ws.Range("M135").Value = -7087.9998
$ws.Range("H136").Value = 130780
$ws.Range("J136").Value = 130780
$ws.Range("L136").Value = 130780
$ws.Range("N136").Value = -140980
$ws.Range("H137").Value = 3742.25
$ws.Range("I137").Value = 1220.9524
$ws.Range("J137").Value = 8555.637000000001
$ws.Range("K137").Value = 3662.857199999999
$ws.Range("L137").Value = 25666.911
$ws.Range("M137").Value = -1112.857199999999
$ws.Range("N137").Value = -30766.911

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15386345
$ws.Range("I32").Value = 18183180
$ws.Range("K32").Value = 18183180
$ws.Range("M32").Value = -18182893
$ws.Range("H37").Value = 22666
$ws.Range("I37").Value = 22666
$ws.Range("K37").Value = 22666
$ws.Range("M37").Value = -22393
$ws.Range("H74").Value = 333707740
$ws.Range("J74").Value = 3125
$ws.Range("L74").Value = 3125
$ws.Range("N74").Value = -4873
$ws.Range("H77").Value = 333707740
$ws.Range("J77").Value = 3125
$ws.Range("L77").Value = 15625
$ws.Range("N77").Value = -24361
$ws.Range("H110").Value = 24467.777
$ws.Range("I110").Value = 33491
$ws.Range("J110").Value = 6421.3335
$ws.Range("K110").Value = 33491
$ws.Range("L110").Value = 6421.3335
$ws.Range("M110").Value = -31446
$ws.Range("N110").Value = -10511.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1326
$ws.Range("I12").Value = 334.8
$ws.Range("J12").Value = 2317.2
$ws.Range("K12").Value = 334.8
$ws.Range("L12").Value = 2317.2
$ws.Range("M12").Value = -166.8
$ws.Range("N12").Value = -2653.2
$ws.Range("H22").Value = 206.57143
$ws.Range("I22").Value = 206.57143
$ws.Range("K22").Value = 206.57143
$ws.Range("M22").Value = -33.57142999999999
$ws.Range("H105").Value = 12978.333
$ws.Range("I105").Value = 18517.5
$ws.Range("K105").Value = 18517.5
$ws.Range("M105").Value = -16770.5
$ws.Range("H132").Value = 101441.664
$ws.Range("J132").Value = 101441.664
$ws.Range("L132").Value = 101441.664
$ws.Range("N132").Value = -111561.664
$ws.Range("H134").Value = 2477.0476
$ws.Range("I134").Value = 2460.9
$ws.Range("K134").Value = 7382.700000000001
$ws.Range("M134").Value = -4847.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 12966.875
$ws.Range("I22").Value = 20360.2
$ws.Range("J22").Value = 644.6667
$ws.Range("K22").Value = 20360.2
$ws.Range("L22").Value = 644.6667
$ws.Range("M22").Value = -20010.2
$ws.Range("N22").Value = -1344.6667
$ws.Range("H58").Value = 3700.2144
$ws.Range("I58").Value = 4491.4443
$ws.Range("K58").Value = 4491.4443
$ws.Range("M58").Value = -4288.4443
$ws.Range("H132").Value = 2137.7715
$ws.Range("I132").Value = 1949.129
$ws.Range("J132").Value = 3599.75
$ws.Range("K132").Value = 5847.387
$ws.Range("L132").Value = 10799.25
$ws.Range("M132").Value = -3317.387
$ws.Range("N132").Value = -15859.25
$ws.Range("H133").Value = 50217
$ws.Range("J133").Value = 60325.5
$ws.Range("L133").Value = 60325.5
$ws.Range("N133").Value = -65385.5
$ws.Range("H134").Value = 7957.7144
$ws.Range("I134").Value = 7341.4
$ws.Range("K134").Value = 22024.2
$ws.Range("M134").Value = -19489.2
$ws.Range("H136").Value = 3700.2144
$ws.Range("I136").Value = 4491.4443
$ws.Range("K136").Value = 13474.3329
$ws.Range("M136").Value = -10924.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1406.25
$ws.Range("I97").Value = 1406.25
$ws.Range("K97").Value = 1406.25
$ws.Range("M97").Value = -910.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 67642.5
$ws.Range("I109").Value = 80000
$ws.Range("J109").Value = 55285
$ws.Range("K109").Value = 80000
$ws.Range("L109").Value = 55285
$ws.Range("M109").Value = -78613
$ws.Range("N109").Value = -58059
$ws.Range("H132").Value = 66670544
$ws.Range("I132").Value = 4198.273
$ws.Range("J132").Value = 250003000
$ws.Range("K132").Value = 12594.819
$ws.Range("L132").Value = 750009000
$ws.Range("M132").Value = -10064.819
$ws.Range("N132").Value = -750014060
$ws.Range("H136").Value = 911554.0600000001
$ws.Range("I136").Value = 1430535.4
$ws.Range("J136").Value = 3336.75
$ws.Range("K136").Value = 4291606.199999999
$ws.Range("L136").Value = 10010.25
$ws.Range("M136").Value = -4289056.199999999
$ws.Range("N136").Value = -15110.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 13000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H62").Value = 18900
$ws.Range("J62").Value = 18900
$ws.Range("L62").Value = 18900
$ws.Range("N62").Value = -20148
$ws.Range("H65").Value = 18900
$ws.Range("J65").Value = 18900
$ws.Range("L65").Value = 94500
$ws.Range("N65").Value = -100740
$ws.Range("H81").Value = 3861.6072
$ws.Range("J81").Value = 7249.875
$ws.Range("L81").Value = 14499.75
$ws.Range("N81").Value = -16621.75
$ws.Range("H84").Value = 3861.6072
$ws.Range("J84").Value = 7249.875
$ws.Range("L84").Value = 72498.75
$ws.Range("N84").Value = -83106.75
$ws.Range("H132").Value = 1856.1538
$ws.Range("I132").Value = 1802.5
$ws.Range("K132").Value = 5407.5
$ws.Range("M132").Value = -2877.5
$ws.Range("H136").Value = 1699.4348
$ws.Range("I136").Value = 1499.75
$ws.Range("J136").Value = 2155.8572
$ws.Range("K136").Value = 4499.25
$ws.Range("L136").Value = 6467.571599999999
$ws.Range("M136").Value = -1949.25
$ws.Range("N136").Value = -11567.5716

Write-Host "Applied Famfrit_Profits updates"